# New weekly observation for "Orégano" (Mercado Mayorista Lo Valledor de
# Santiago) inserted as row 307, pushing the existing rows 307-321 down to
# 308-322 (dimension grows from A1:R321 to A1:R322).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 307 - everything below shifts down.
$ws.Rows.Item(307).Insert()

$ws.Cells.Item(307, 1).Value  = 6
$ws.Cells.Item(307, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(307, 3).Value  = "Metropolitana"
$ws.Cells.Item(307, 4).Value  = 45041
$ws.Cells.Item(307, 5).Value  = 13
$ws.Cells.Item(307, 6).Value  = 100112029
$ws.Cells.Item(307, 7).Value  = "Orégano"
$ws.Cells.Item(307, 8).Value  = "Sin especificar"
$ws.Cells.Item(307, 9).Value  = "Primera"
$ws.Cells.Item(307, 10).Value = 45
$ws.Cells.Item(307, 11).Value = 16000
$ws.Cells.Item(307, 12).Value = 17000
$ws.Cells.Item(307, 13).Value = 16467
$ws.Cells.Item(307, 14).Value = "$/docena de atados"
$ws.Cells.Item(307, 15).Value = "Región Metropolitana"
$ws.Cells.Item(307, 16).Value = 5489
$ws.Cells.Item(307, 17).Value = 3
$ws.Cells.Item(307, 18).Value = "Hortaliza"
